# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Update Rule column value on row 11 from "R40" to "1".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "'1"
